# Ra_Stock_5 concentration correction: rows 57-81 used the wrong stock-5
# radium concentration; recompute Cs/sCs/TotAct/fSorb/sfSorb (E:I) and bump
# the Date of Input (L) to reflect the corrected re-processing date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57
$ws.Range("E57").Value = 270.81622261908865
$ws.Range("F57").Value = 3.5223545816029764
$ws.Range("G57").Value = 9.778649999999999
$ws.Range("H57").Value = 0.84179714838932884
$ws.Range("I57").Value = 0.0029176117819023442
$ws.Range("L57").Value = 42844

# Row 58
$ws.Range("E58").Value = 1543.9279231594101
$ws.Range("F58").Value = 27.390610686569982
$ws.Range("G58").Value = 48.963600000000007
$ws.Range("H58").Value = 0.94362105126578921
$ws.Range("I58").Value = 0.0059433182421027844
$ws.Range("L58").Value = 42844

# Row 59
$ws.Range("E59").Value = 3025.0443274560475
$ws.Range("F59").Value = 50.819988198928094
$ws.Range("G59").Value = 97.786500000000004
$ws.Range("H59").Value = 0.93519060357189066
$ws.Range("I59").Value = 0.010052073987640393
$ws.Range("L59").Value = 42844

# Row 60
$ws.Range("E60").Value = 7332.0338173200971
$ws.Range("F60").Value = 86.333322355887105
$ws.Range("G60").Value = 244.81799999999998
$ws.Range("H60").Value = 0.90146617923715444
$ws.Range("I60").Value = 0.015231309369218682
$ws.Range("L60").Value = 42844

# Row 61
$ws.Range("E61").Value = 14679.955354889733
$ws.Range("F61").Value = 130.33299320161686
$ws.Range("G61").Value = 489.63599999999997
$ws.Range("H61").Value = 0.9003949709257405
$ws.Range("I61").Value = 0.0058232625164407106
$ws.Range("L61").Value = 42844

# Row 62
$ws.Range("E62").Value = 92.47011910287894
$ws.Range("F62").Value = 4.7370138367286376
$ws.Range("G62").Value = 9.778649999999999
$ws.Range("H62").Value = 0.38149406496338573
$ws.Range("I62").Value = 0.042039687351188987
$ws.Range("L62").Value = 42844

# Row 63
$ws.Range("E63").Value = 627.08534477244473
$ws.Range("F63").Value = 32.954498809913645
$ws.Range("G63").Value = 48.963600000000007
$ws.Range("H63").Value = 0.54079969318947774
$ws.Range("I63").Value = 0.011307961032332713
$ws.Range("L63").Value = 42844

# Row 64
$ws.Range("E64").Value = 1334.6731089275274
$ws.Range("F64").Value = 206.78881287892077
$ws.Range("G64").Value = 97.786500000000004
$ws.Range("H64").Value = 0.54114371551005247
$ws.Range("I64").Value = 0.0044944264440218446
$ws.Range("L64").Value = 42844

# Row 65
$ws.Range("E65").Value = 3279.9381783940717
$ws.Range("F65").Value = 576.72918203223719
$ws.Range("G65").Value = 244.81799999999998
$ws.Range("H65").Value = 0.54503464293473003
$ws.Range("I65").Value = 0.02062412086264484
$ws.Range("L65").Value = 42844

# Row 66
$ws.Range("E66").Value = 6646.1784208833124
$ws.Range("F66").Value = 1187.0797757543103
$ws.Range("G66").Value = 489.63599999999997
$ws.Range("H66").Value = 0.53893602092389647
$ws.Range("I66").Value = 0.0017379476258423448
$ws.Range("L66").Value = 42844

# Row 67
$ws.Range("E67").Value = 0.11444991418286077
$ws.Range("F67").Value = 0.01070047154163809
$ws.Range("G67").Value = 9.778649999999999
$ws.Range("H67").Value = 0.43834237035580376
$ws.Range("I67").Value = 0.017670864210681043
$ws.Range("L67").Value = 42844

# Row 68
$ws.Range("E68").Value = 0.73587745446631947
$ws.Range("F68").Value = 0.025045985577600658
$ws.Range("G68").Value = 48.963600000000007
$ws.Range("H68").Value = 0.54071101632868901
$ws.Range("I68").Value = 0.0047569624366485246
$ws.Range("L68").Value = 42844

# Row 69
$ws.Range("E69").Value = 1.3242049550090049
$ws.Range("F69").Value = 0.22627993517601505
$ws.Range("G69").Value = 97.786500000000004
$ws.Range("H69").Value = 0.5542746649141852
$ws.Range("I69").Value = 0.0062283865550030989
$ws.Range("L69").Value = 42844

# Row 70
$ws.Range("E70").Value = 3.7933390740114326
$ws.Range("F70").Value = 0.4369919422995962
$ws.Range("G70").Value = 244.81799999999998
$ws.Range("H70").Value = 0.55710007035718179
$ws.Range("I70").Value = 0.0068005762380074983
$ws.Range("L70").Value = 42844

# Row 71
$ws.Range("E71").Value = 7.1863934914274834
$ws.Range("F71").Value = 1.4261909061889864
$ws.Range("G71").Value = 489.63599999999997
$ws.Range("H71").Value = 0.5764148930595675
$ws.Range("I71").Value = 0.0031099084834229291
$ws.Range("L71").Value = 42844

# Row 72
$ws.Range("E72").Value = 120.2965404934521
$ws.Range("F72").Value = 22.847191531433218
$ws.Range("G72").Value = 9.778649999999999
$ws.Range("H72").Value = 0.469249011180983
$ws.Range("I72").Value = 0.022708007255614518
$ws.Range("L72").Value = 42844

# Row 73
$ws.Range("E73").Value = 772.37480257413097
$ws.Range("F73").Value = 51.904722149829389
$ws.Range("G73").Value = 48.963600000000007
$ws.Range("H73").Value = 0.57150875879833707
$ws.Range("I73").Value = 0.0038679315018278803
$ws.Range("L73").Value = 42844

# Row 74
$ws.Range("E74").Value = 1557.2607790595678
$ws.Range("F74").Value = 191.88410186275289
$ws.Range("G74").Value = 97.786500000000004
$ws.Range("H74").Value = 0.58323433537480052
$ws.Range("I74").Value = 0.0053911615325178009
$ws.Range("L74").Value = 42844

# Row 75
$ws.Range("E75").Value = 4122.4980908609041
$ws.Range("F75").Value = 828.50446916051305
$ws.Range("G75").Value = 244.81799999999998
$ws.Range("H75").Value = 0.59487480074657484
$ws.Range("I75").Value = 0.021472861288731531
$ws.Range("L75").Value = 42844

# Row 76
$ws.Range("E76").Value = 7797.2186700330813
$ws.Range("F76").Value = 768.60264682482546
$ws.Range("G76").Value = 489.63599999999997
$ws.Range("H76").Value = 0.59614898791591653
$ws.Range("I76").Value = 0.0031328926411537312
$ws.Range("L76").Value = 42844

# Row 77
$ws.Range("E77").Value = 127.81352619499012
$ws.Range("F77").Value = 17.864401649697054
$ws.Range("G77").Value = 9.778649999999999
$ws.Range("H77").Value = 0.5178425389738156
$ws.Range("I77").Value = 0.036599482021059691
$ws.Range("L77").Value = 42844

# Row 78
$ws.Range("E78").Value = 727.21203336213875
$ws.Range("F78").Value = 85.942098216403082
$ws.Range("G78").Value = 48.963600000000007
$ws.Range("H78").Value = 0.62208412840144678
$ws.Range("I78").Value = 0.013625560510124042
$ws.Range("L78").Value = 42844

# Row 79
$ws.Range("E79").Value = 1477.9033340471351
$ws.Range("F79").Value = 97.898483982226224
$ws.Range("G79").Value = 97.786500000000004
$ws.Range("H79").Value = 0.63300654045813587
$ws.Range("I79").Value = 0.0015685037377726563
$ws.Range("L79").Value = 42844

# Row 80
$ws.Range("E80").Value = 3755.1320098543606
$ws.Range("F80").Value = 1084.4336807657073
$ws.Range("G80").Value = 244.81799999999998
$ws.Range("H80").Value = 0.60939630849511151
$ws.Range("I80").Value = 0.018529052300006489
$ws.Range("L80").Value = 42844

# Row 81
$ws.Range("E81").Value = 7016.7029571436624
$ws.Range("F81").Value = 427.21305716109964
$ws.Range("G81").Value = 489.63599999999997
$ws.Range("H81").Value = 0.60923084026543728
$ws.Range("I81").Value = 0.015296611248176354
$ws.Range("L81").Value = 42844

# Match the author's final selection/view state on the corrected Date column
$ws.Range("L57:L81").Select() | Out-Null
